$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from 46075 (2026-02-22)
# to 46076 (2026-02-23) for every data row (rows 2 through 23).
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 46076
}
